$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell (never itself edited by this script) for restoring the default,
# unstyled format after forcing text entry on numeric-looking price strings,
# so cells that need the "@" text-format trick keep their original style index.
$donorStyle = $ws.Range("C2").Style

# --- Cells whose new value is a plain text / multi-dot string (no locale-number risk) ---
$ws.Range("D2").Value = "61.073.91"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "3.372.22"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("E6").Value = "  +9.58%  "
$ws.Range("E7").Value = "  +2.58%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +4.85%  "
$ws.Range("E10").Value = "  +8.90%  "
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "3.883.60"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "3.359.11"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "61.007.68"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("E20").Value = "  +8.50%  "
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("E22").Value = "  +10.53%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +12.03%  "
$ws.Range("E27").Value = "  +8.79%  "
$ws.Range("E28").Value = "  -4.16%  "
$ws.Range("E29").Value = "  -7.58%  "
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -1.00%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("E40").Value = "  -3.96%  "
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E48").Value = "  -2.38%  "
$ws.Range("D49").Value = "2.122.50"
$ws.Range("E49").Value = "  -3.56%  "
$ws.Range("E51").Value = "  +0.13%  "

# --- Price cells whose new value parses as a plain number; force Text format first
#     so Excel keeps the original formatted-string representation (e.g. "8.33"),
#     then restore the original (unstyled) cell style. ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "406.16"
$ws.Range("D5").Style = $donorStyle
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.74"
$ws.Range("D6").Style = $donorStyle
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.591"
$ws.Range("D7").Style = $donorStyle
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.672"
$ws.Range("D9").Style = $donorStyle
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("D10").Style = $donorStyle
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.45"
$ws.Range("D11").Style = $donorStyle
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.33"
$ws.Range("D14").Style = $donorStyle
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.62"
$ws.Range("D15").Style = $donorStyle
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.09"
$ws.Range("D19").Style = $donorStyle
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000128"
$ws.Range("D20").Style = $donorStyle
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.22"
$ws.Range("D21").Style = $donorStyle
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "84.05"
$ws.Range("D22").Style = $donorStyle
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "307.81"
$ws.Range("D23").Style = $donorStyle
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.72"
$ws.Range("D24").Style = $donorStyle
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.13"
$ws.Range("D25").Style = $donorStyle
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.78"
$ws.Range("D26").Style = $donorStyle
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.43"
$ws.Range("D28").Style = $donorStyle
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.54"
$ws.Range("D29").Style = $donorStyle
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.28"
$ws.Range("D33").Style = $donorStyle
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.23"
$ws.Range("D34").Style = $donorStyle
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("D35").Style = $donorStyle
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0480"
$ws.Range("D36").Style = $donorStyle
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.81"
$ws.Range("D37").Style = $donorStyle
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.997"
$ws.Range("D38").Style = $donorStyle
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("D39").Style = $donorStyle
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.90"
$ws.Range("D40").Style = $donorStyle
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.99"
$ws.Range("D41").Style = $donorStyle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "137.12"
$ws.Range("D42").Style = $donorStyle
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.02"
$ws.Range("D44").Style = $donorStyle
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.286"
$ws.Range("D45").Style = $donorStyle
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.65"
$ws.Range("D46").Style = $donorStyle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("D47").Style = $donorStyle
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.40"
$ws.Range("D48").Style = $donorStyle
